$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 9,28
$arr[0,0] = "A1"
$arr[0,1] = "M3"
$arr[0,2] = "DO"
$arr[0,3] = "A1"
$arr[0,4] = "M1"
$arr[0,5] = "M1"
$arr[0,6] = "M1"
$arr[0,7] = "A1"
$arr[0,8] = "M1"
$arr[0,9] = "A1"
$arr[0,10] = "A1"
$arr[0,11] = "A1"
$arr[0,12] = "M3"
$arr[0,13] = "DO"
$arr[0,14] = "M3"
$arr[0,15] = "A1"
$arr[0,16] = "A1"
$arr[0,17] = "A1"
$arr[0,18] = "A1"
$arr[0,19] = "M1"
$arr[0,20] = "DO"
$arr[0,21] = "A1"
$arr[0,22] = "DO"
$arr[0,23] = "M3"
$arr[0,24] = "A1"
$arr[0,25] = "A1"
$arr[0,26] = "M1"
$arr[0,27] = "A1"
$arr[1,0] = "A2"
$arr[1,1] = "M2"
$arr[1,2] = "A1"
$arr[1,3] = "M2"
$arr[1,4] = "DO"
$arr[1,5] = "M1"
$arr[1,6] = "M3"
$arr[1,7] = "A2"
$arr[1,8] = "M3"
$arr[1,9] = "M2"
$arr[1,10] = "DO"
$arr[1,11] = "M1"
$arr[1,12] = "M1"
$arr[1,13] = "M1"
$arr[1,14] = "M1"
$arr[1,15] = "M3"
$arr[1,16] = "DO"
$arr[1,17] = "M2"
$arr[1,18] = "M2"
$arr[1,19] = "A1"
$arr[1,20] = "M1"
$arr[1,21] = "A2"
$arr[1,22] = "M3"
$arr[1,23] = "DO"
$arr[1,24] = "M2"
$arr[1,25] = "M1"
$arr[1,26] = "M2"
$arr[1,27] = "M2"
$arr[2,0] = "M1"
$arr[2,1] = "DO"
$arr[2,2] = "M1"
$arr[2,3] = "M3"
$arr[2,4] = "M1"
$arr[2,5] = "A1"
$arr[2,6] = "M1"
$arr[2,7] = "M3"
$arr[2,8] = "A1"
$arr[2,9] = "DO"
$arr[2,10] = "M1"
$arr[2,11] = "M1"
$arr[2,12] = "M1"
$arr[2,13] = "M1"
$arr[2,14] = "M2"
$arr[2,15] = "M1"
$arr[2,16] = "M1"
$arr[2,17] = "M1"
$arr[2,18] = "DO"
$arr[2,19] = "A1"
$arr[2,20] = "M3"
$arr[2,21] = "M1"
$arr[2,22] = "DO"
$arr[2,23] = "M1"
$arr[2,24] = "M1"
$arr[2,25] = "M1"
$arr[2,26] = "A1"
$arr[2,27] = "M3"
$arr[3,0] = "M1"
$arr[3,1] = "DO"
$arr[3,2] = "M1"
$arr[3,3] = "M1"
$arr[3,4] = "M3"
$arr[3,5] = "M2"
$arr[3,6] = "A2"
$arr[3,7] = "M1"
$arr[3,8] = "DO"
$arr[3,9] = "M1"
$arr[3,10] = "M1"
$arr[3,11] = "M3"
$arr[3,12] = "M3"
$arr[3,13] = "A2"
$arr[3,14] = "DO"
$arr[3,15] = "M1"
$arr[3,16] = "M1"
$arr[3,17] = "M2"
$arr[3,18] = "M1"
$arr[3,19] = "M3"
$arr[3,20] = "A2"
$arr[3,21] = "M1"
$arr[3,22] = "M2"
$arr[3,23] = "DO"
$arr[3,24] = "M2"
$arr[3,25] = "M3"
$arr[3,26] = "M2"
$arr[3,27] = "M3"
$arr[4,0] = "M1"
$arr[4,1] = "M1"
$arr[4,2] = "DO"
$arr[4,3] = "M1"
$arr[4,4] = "A1"
$arr[4,5] = "M3"
$arr[4,6] = "M1"
$arr[4,7] = "M3"
$arr[4,8] = "M1"
$arr[4,9] = "M1"
$arr[4,10] = "M1"
$arr[4,11] = "DO"
$arr[4,12] = "M2"
$arr[4,13] = "M1"
$arr[4,14] = "A2"
$arr[4,15] = "DO"
$arr[4,16] = "M3"
$arr[4,17] = "M1"
$arr[4,18] = "M1"
$arr[4,19] = "M1"
$arr[4,20] = "M1"
$arr[4,21] = "A2"
$arr[4,22] = "A2"
$arr[4,23] = "A2"
$arr[4,24] = "DO"
$arr[4,25] = "M1"
$arr[4,26] = "M3"
$arr[4,27] = "M1"
$arr[5,0] = "DO"
$arr[5,1] = "A1"
$arr[5,2] = "M3"
$arr[5,3] = "A2"
$arr[5,4] = "A1"
$arr[5,5] = "A2"
$arr[5,6] = "A2"
$arr[5,7] = "M3"
$arr[5,8] = "A1"
$arr[5,9] = "DO"
$arr[5,10] = "A2"
$arr[5,11] = "A2"
$arr[5,12] = "A2"
$arr[5,13] = "A2"
$arr[5,14] = "A1"
$arr[5,15] = "DO"
$arr[5,16] = "M3"
$arr[5,17] = "A1"
$arr[5,18] = "A1"
$arr[5,19] = "M1"
$arr[5,20] = "A1"
$arr[5,21] = "DO"
$arr[5,22] = "M3"
$arr[5,23] = "A1"
$arr[5,24] = "A1"
$arr[5,25] = "A1"
$arr[5,26] = "A1"
$arr[5,27] = "A1"
$arr[6,0] = "A1"
$arr[6,1] = "A2"
$arr[6,2] = "A1"
$arr[6,3] = "A1"
$arr[6,4] = "A2"
$arr[6,5] = "M3"
$arr[6,6] = "DO"
$arr[6,7] = "A1"
$arr[6,8] = "A1"
$arr[6,9] = "A1"
$arr[6,10] = "DO"
$arr[6,11] = "A1"
$arr[6,12] = "M3"
$arr[6,13] = "M3"
$arr[6,14] = "A2"
$arr[6,15] = "M3"
$arr[6,16] = "A2"
$arr[6,17] = "A1"
$arr[6,18] = "A2"
$arr[6,19] = "M2"
$arr[6,20] = "DO"
$arr[6,21] = "M1"
$arr[6,22] = "A2"
$arr[6,23] = "A2"
$arr[6,24] = "A1"
$arr[6,25] = "A2"
$arr[6,26] = "M3"
$arr[6,27] = "DO"
$arr[7,0] = "DO"
$arr[7,1] = "M3"
$arr[7,2] = "M2"
$arr[7,3] = "M2"
$arr[7,4] = "M2"
$arr[7,5] = "A2"
$arr[7,6] = "A2"
$arr[7,7] = "A1"
$arr[7,8] = "DO"
$arr[7,9] = "M1"
$arr[7,10] = "M3"
$arr[7,11] = "M2"
$arr[7,12] = "A1"
$arr[7,13] = "A2"
$arr[7,14] = "DO"
$arr[7,15] = "M1"
$arr[7,16] = "M3"
$arr[7,17] = "M2"
$arr[7,18] = "M2"
$arr[7,19] = "A2"
$arr[7,20] = "A1"
$arr[7,21] = "DO"
$arr[7,22] = "M3"
$arr[7,23] = "M1"
$arr[7,24] = "M2"
$arr[7,25] = "M1"
$arr[7,26] = "A1"
$arr[7,27] = "A1"
$arr[8,0] = "M1"
$arr[8,1] = "M2"
$arr[8,2] = "M3"
$arr[8,3] = "M1"
$arr[8,4] = "DO"
$arr[8,5] = "M2"
$arr[8,6] = "M2"
$arr[8,7] = "DO"
$arr[8,8] = "M1"
$arr[8,9] = "A1"
$arr[8,10] = "A1"
$arr[8,11] = "M2"
$arr[8,12] = "A1"
$arr[8,13] = "M3"
$arr[8,14] = "M1"
$arr[8,15] = "A1"
$arr[8,16] = "A1"
$arr[8,17] = "DO"
$arr[8,18] = "M3"
$arr[8,19] = "M2"
$arr[8,20] = "M2"
$arr[8,21] = "M2"
$arr[8,22] = "M1"
$arr[8,23] = "M2"
$arr[8,24] = "DO"
$arr[8,25] = "M3"
$arr[8,26] = "M3"
$arr[8,27] = "M1"
$ws.Range("B2:AC10").Value = $arr
